{"js": "// Update the date and the 25 multiplication problems in the table.\n// Some \"before\" values are duplicated in the document (e.g. \"55\u00d785=4675\"\n// appears twice) but map to different \"after\" values, so replacements are\n// applied in document order using the index of each duplicate occurrence.\n\nconst replacements = [\n  { find: \"2024-06-14 Friday\", replace: \"2024-06-15 Saturday\", occurrence: 0 },\n  { find: \"64\u00d761=3904\", replace: \"52\u00d764=3328\", occurrence: 0 },\n  { find: \"31\u00d728=868\", replace: \"98\u00d792=9016\", occurrence: 0 },\n  { find: \"70\u00d772=5040\", replace: \"18\u00d712=216\", occurrence: 0 },\n  { find: \"50\u00d798=4900\", replace: \"71\u00d774=5254\", occurrence: 0 },\n  { find: \"88\u00d791=8008\", replace: \"64\u00d765=4160\", occurrence: 0 },\n  { find: \"97\u00d740=3880\", replace: \"61\u00d728=1708\", occurrence: 0 },\n  { find: \"92\u00d768=6256\", replace: \"32\u00d752=1664\", occurrence: 0 },\n  { find: \"69\u00d739=2691\", replace: \"63\u00d729=1827\", occurrence: 0 },\n  { find: \"55\u00d785=4675\", replace: \"36\u00d738=1368\", occurrence: 0 },\n  { find: \"56\u00d787=4872\", replace: \"75\u00d716=1200\", occurrence: 0 },\n  { find: \"20\u00d797=1940\", replace: \"94\u00d756=5264\", occurrence: 0 },\n  { find: \"83\u00d783=6889\", replace: \"78\u00d770=5460\", occurrence: 0 },\n  { find: \"80\u00d718=1440\", replace: \"88\u00d711=968\", occurrence: 0 },\n  { find: \"85\u00d745=3825\", replace: \"13\u00d798=1274\", occurrence: 0 },\n  { find: \"75\u00d722=1650\", replace: \"55\u00d799=5445\", occurrence: 0 },\n  { find: \"73\u00d784=6132\", replace: \"25\u00d751=1275\", occurrence: 0 },\n  { find: \"74\u00d748=3552\", replace: \"70\u00d740=2800\", occurrence: 0 },\n  { find: \"18\u00d775=1350\", replace: \"80\u00d774=5920\", occurrence: 0 },\n  { find: \"65\u00d757=3705\", replace: \"73\u00d739=2847\", occurrence: 0 },\n  { find: \"43\u00d742=1806\", replace: \"13\u00d725=325\", occurrence: 0 },\n  { find: \"55\u00d785=4675\", replace: \"60\u00d737=2220\", occurrence: 1 },\n  { find: \"50\u00d794=4700\", replace: \"71\u00d768=4828\", occurrence: 0 },\n  { find: \"12\u00d726=312\", replace: \"93\u00d771=6603\", occurrence: 0 },\n  { find: \"84\u00d742=3528\", replace: \"88\u00d767=5896\", occurrence: 0 },\n  { find: \"40\u00d724=960\", replace: \"82\u00d722=1804\", occurrence: 0 },\n];\n\nconst body = context.document.body;\n\n// Group by \"find\" text so duplicate searches are issued only once and\n// resolved through the occurrence index of the search results collection\n// (which is returned in document order).\nconst byFind = new Map();\nfor (const r of replacements) {\n  if (!byFind.has(r.find)) byFind.set(r.find, []);\n  byFind.get(r.find).push(r);\n}\n\nconst searchResultsByFind = new Map();\nfor (const find of byFind.keys()) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  searchResultsByFind.set(find, results);\n}\n\nawait context.sync();\n\nfor (const [find, entries] of byFind.entries()) {\n  const results = searchResultsByFind.get(find);\n  for (const entry of entries) {\n    const range = results.items[entry.occurrence];\n    range.insertText(entry.replace, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and the 25 multiplication problems in the table.\n# Note: \"55\u00d785=4675\" occurs twice in the document (different cells) and\n# maps to two different replacement values, so we issue two sequential\n# single-replace Find/Execute calls for that phrase; because the Find\n# range cursor advances past each match, the second call automatically\n# lands on the next (second) occurrence in document order.\n\n$d = $word.ActiveDocument\n\nfunction Replace-OneOccurrence($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 1) | Out-Null\n}\n\nReplace-OneOccurrence \"2024-06-14 Friday\" \"2024-06-15 Saturday\"\n\nReplace-OneOccurrence \"64\u00d761=3904\" \"52\u00d764=3328\"\nReplace-OneOccurrence \"31\u00d728=868\" \"98\u00d792=9016\"\nReplace-OneOccurrence \"70\u00d772=5040\" \"18\u00d712=216\"\nReplace-OneOccurrence \"50\u00d798=4900\" \"71\u00d774=5254\"\nReplace-OneOccurrence \"88\u00d791=8008\" \"64\u00d765=4160\"\n\nReplace-OneOccurrence \"97\u00d740=3880\" \"61\u00d728=1708\"\nReplace-OneOccurrence \"92\u00d768=6256\" \"32\u00d752=1664\"\nReplace-OneOccurrence \"69\u00d739=2691\" \"63\u00d729=1827\"\nReplace-OneOccurrence \"55\u00d785=4675\" \"36\u00d738=1368\"\nReplace-OneOccurrence \"56\u00d787=4872\" \"75\u00d716=1200\"\n\nReplace-OneOccurrence \"20\u00d797=1940\" \"94\u00d756=5264\"\nReplace-OneOccurrence \"83\u00d783=6889\" \"78\u00d770=5460\"\nReplace-OneOccurrence \"80\u00d718=1440\" \"88\u00d711=968\"\nReplace-OneOccurrence \"85\u00d745=3825\" \"13\u00d798=1274\"\nReplace-OneOccurrence \"75\u00d722=1650\" \"55\u00d799=5445\"\n\nReplace-OneOccurrence \"73\u00d784=6132\" \"25\u00d751=1275\"\nReplace-OneOccurrence \"74\u00d748=3552\" \"70\u00d740=2800\"\nReplace-OneOccurrence \"18\u00d775=1350\" \"80\u00d774=5920\"\nReplace-OneOccurrence \"65\u00d757=3705\" \"73\u00d739=2847\"\nReplace-OneOccurrence \"43\u00d742=1806\" \"13\u00d725=325\"\n\nReplace-OneOccurrence \"55\u00d785=4675\" \"60\u00d737=2220\"\nReplace-OneOccurrence \"50\u00d794=4700\" \"71\u00d768=4828\"\nReplace-OneOccurrence \"12\u00d726=312\" \"93\u00d771=6603\"\nReplace-OneOccurrence \"84\u00d742=3528\" \"88\u00d767=5896\"\nReplace-OneOccurrence \"40\u00d724=960\" \"82\u00d722=1804\"\n"}
